# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    2  = 855
    4  = 2168
    6  = 12604
    9  = 505
    10 = 461
    11 = 1150
    12 = 949
    13 = 13660
    14 = 13983
    19 = 13
    23 = 1053
    26 = 621
    27 = 5100
    29 = 253
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
